$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before current row 115 (shifting old rows 115-180 down to 117-182)
$ws.Rows("115:116").Insert()

# New row 115: Camote, 1a (guarda), Región del Maule
$ws.Cells.Item(115, 1).Value2  = 7
$ws.Cells.Item(115, 2).Value2  = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(115, 3).Value2  = "Ñuble"
$ws.Cells.Item(115, 4).Value2  = 44813
$ws.Cells.Item(115, 5).Value2  = 16
$ws.Cells.Item(115, 6).Value2  = 100112045
$ws.Cells.Item(115, 7).Value2  = "Zapallo"
$ws.Cells.Item(115, 8).Value2  = "Camote"
$ws.Cells.Item(115, 9).Value2  = "1a (guarda)"
$ws.Cells.Item(115, 10).Value2 = 300
$ws.Cells.Item(115, 11).Value2 = 800
$ws.Cells.Item(115, 12).Value2 = 900
$ws.Cells.Item(115, 13).Value2 = 850
$ws.Cells.Item(115, 14).Value2 = "`$/kilo (volumen en unidades)"
$ws.Cells.Item(115, 15).Value2 = "Región del Maule"
$ws.Cells.Item(115, 16).Value2 = 850
$ws.Cells.Item(115, 17).Value2 = 1
$ws.Cells.Item(115, 18).Value2 = "Hortaliza"

# New row 116: Camote, 2a (guarda), Región del Maule
$ws.Cells.Item(116, 1).Value2  = 7
$ws.Cells.Item(116, 2).Value2  = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(116, 3).Value2  = "Ñuble"
$ws.Cells.Item(116, 4).Value2  = 44813
$ws.Cells.Item(116, 5).Value2  = 16
$ws.Cells.Item(116, 6).Value2  = 100112045
$ws.Cells.Item(116, 7).Value2  = "Zapallo"
$ws.Cells.Item(116, 8).Value2  = "Camote"
$ws.Cells.Item(116, 9).Value2  = "2a (guarda)"
$ws.Cells.Item(116, 10).Value2 = 200
$ws.Cells.Item(116, 11).Value2 = 700
$ws.Cells.Item(116, 12).Value2 = 700
$ws.Cells.Item(116, 13).Value2 = 700
$ws.Cells.Item(116, 14).Value2 = "`$/kilo (volumen en unidades)"
$ws.Cells.Item(116, 15).Value2 = "Región del Maule"
$ws.Cells.Item(116, 16).Value2 = 700
$ws.Cells.Item(116, 17).Value2 = 1
$ws.Cells.Item(116, 18).Value2 = "Hortaliza"
